# Generate Report for Handoff
# Adds two new "Ready for handoff" rows (a .png handback and an .md handback)
# to the Overview / zh-cn / de-de sheets, each with the associated hyperlinks,
# and refreshes the existing first data row to the latest uuid-named source
# file that was (re)handed off in this run.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# ---------------------------------------------------------------------------
# New source file identities handed off in this run.
# ---------------------------------------------------------------------------
$file2 = "1c05f62d-3a8b-4fb7-a33d-6b94be88534d.png"
$file3 = "373f3758-a65d-48f4-b48c-9a84c004ef77.png"
$file4 = "f5af7ab7-db05-437c-a321-55e479a000fe.md"

$zhTarget2 = "e48444fb348d0a76e1bb1f1e63ba08f40b259613.png"
$zhTarget3 = "32e97e2c7ae78dbc312d17ef5ac7fbf7ffa14d3c.png"
$zhTarget4 = "f5af7ab7-db05-437c-a321-55e479a000fe.16872519c2f7bf7e9f93fe198086f63062a210c1.zh-cn.xlf"

$deTarget2 = "e48444fb348d0a76e1bb1f1e63ba08f40b259613.png"
$deTarget3 = "32e97e2c7ae78dbc312d17ef5ac7fbf7ffa14d3c.png"
$deTarget4 = "f5af7ab7-db05-437c-a321-55e479a000fe.16872519c2f7bf7e9f93fe198086f63062a210c1.de-de.xlf"

$srcBase   = "https://github.com/OpenLocalizationTest/oltest/blob/f54f15f035d31b4f48676748c5a7f58786432d78/e2e/"
$zhHtBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/33ed0b0ebc163441d513e5f44ea8fc30d326008b/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/"
$deHtBase  = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/a0dc2abadb5722145fa3e18b3cd81a24b27104d5/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/"

$statusReady   = "Ready for handoff"
$handoffTime   = "2016-43-11 22:43:01"
$zhHandoffDt   = "2016-03-11 22:42:58"
$deHandoffDt   = "2016-03-11 22:43:01"
$epoch         = "0001-01-01 00:00:00"

# ---------------------------------------------------------------------------
# Drop every existing hyperlink up front -- this engine's
# `Range(...).Hyperlinks.Delete()` clears the whole sheet's collection
# rather than just the target range, so all links get re-added afterwards
# in final top-to-bottom / left-to-right order.
# ---------------------------------------------------------------------------
$ws1.Range("A1").Hyperlinks.Delete()
$ws2.Range("A1").Hyperlinks.Delete()
$ws3.Range("A1").Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# Sheet 1: Overview
# ---------------------------------------------------------------------------
$ws1.Range("A2").Value = $file2
$ws1.Range("B2").Value = $statusReady
$ws1.Range("C2").Value = $statusReady
$ws1.Range("D2").Value = $handoffTime

$ws1.Range("A3").Value = $file3
$ws1.Range("B3").Value = $statusReady
$ws1.Range("C3").Value = $statusReady
$ws1.Range("D3").Value = $handoffTime

$ws1.Range("A4").Value = $file4
$ws1.Range("B4").Value = $statusReady
$ws1.Range("C4").Value = $statusReady
$ws1.Range("D4").Value = $handoffTime

$ws1.Hyperlinks.Add($ws1.Range("A2"), ($srcBase + $file2), "", "", $file2)
$ws1.Hyperlinks.Add($ws1.Range("A3"), ($srcBase + $file3), "", "", $file3)
$ws1.Hyperlinks.Add($ws1.Range("A4"), ($srcBase + $file4), "", "", $file4)

# ---------------------------------------------------------------------------
# Sheet 2: zh-cn
# ---------------------------------------------------------------------------
$ws2.Range("A2").Value = $file2
$ws2.Range("B2").Value = ".png"
$ws2.Range("C2").Value = $statusReady
$ws2.Range("D2").Value = $zhTarget2
$ws2.Range("E2").Value = $zhHandoffDt
$ws2.Range("H2").Value = $epoch
$ws2.Range("I2").Value = "IsDependency"
$ws2.Range("J2").Value = ("e2e\" + $file4)

$ws2.Range("A3").Value = $file3
$ws2.Range("B3").Value = ".png"
$ws2.Range("C3").Value = $statusReady
$ws2.Range("D3").Value = $zhTarget3
$ws2.Range("E3").Value = $zhHandoffDt
$ws2.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H3").Value = $epoch
$ws2.Range("I3").Value = "IsDependency"
$ws2.Range("J3").Value = ("e2e\" + $file4)

$ws2.Range("A4").Value = $file4
$ws2.Range("B4").Value = ".md"
$ws2.Range("C4").Value = $statusReady
$ws2.Range("D4").Value = $zhTarget4
$ws2.Range("E4").Value = $zhHandoffDt
$ws2.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws2.Range("H4").Value = $epoch
$ws2.Range("I4").Value = "Include"

$ws2.Hyperlinks.Add($ws2.Range("A2"), ($srcBase + $file2), "", "", $file2)
$ws2.Hyperlinks.Add($ws2.Range("B2"), ($srcBase + $file2), "", "", ".png")
$ws2.Hyperlinks.Add($ws2.Range("D2"), ($zhHtBase + $zhTarget2), "", "", $zhTarget2)

$ws2.Hyperlinks.Add($ws2.Range("A3"), ($srcBase + $file3), "", "", $file3)
$ws2.Hyperlinks.Add($ws2.Range("B3"), ($srcBase + $file3), "", "", ".png")
$ws2.Hyperlinks.Add($ws2.Range("D3"), ($zhHtBase + $zhTarget3), "", "", $zhTarget3)

$ws2.Hyperlinks.Add($ws2.Range("A4"), ($srcBase + $file4), "", "", $file4)
$ws2.Hyperlinks.Add($ws2.Range("B4"), ($srcBase + $file4), "", "", ".md")
$ws2.Hyperlinks.Add($ws2.Range("D4"), ($zhHtBase + $zhTarget4), "", "", $zhTarget4)

# ---------------------------------------------------------------------------
# Sheet 3: de-de
# ---------------------------------------------------------------------------
$ws3.Range("A2").Value = $file2
$ws3.Range("B2").Value = ".png"
$ws3.Range("C2").Value = $statusReady
$ws3.Range("D2").Value = $deTarget2
$ws3.Range("E2").Value = $deHandoffDt
$ws3.Range("H2").Value = $epoch
$ws3.Range("I2").Value = "IsDependency"
$ws3.Range("J2").Value = ("e2e\" + $file4)

$ws3.Range("A3").Value = $file3
$ws3.Range("B3").Value = ".png"
$ws3.Range("C3").Value = $statusReady
$ws3.Range("D3").Value = $deTarget3
$ws3.Range("E3").Value = $deHandoffDt
$ws3.Range("E3").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("H3").Value = $epoch
$ws3.Range("I3").Value = "IsDependency"
$ws3.Range("J3").Value = ("e2e\" + $file4)

$ws3.Range("A4").Value = $file4
$ws3.Range("B4").Value = ".md"
$ws3.Range("C4").Value = $statusReady
$ws3.Range("D4").Value = $deTarget4
$ws3.Range("E4").Value = $deHandoffDt
$ws3.Range("E4").NumberFormat = "yyyy-mm-dd HH:mm:ss"
$ws3.Range("H4").Value = $epoch
$ws3.Range("I4").Value = "Include"

$ws3.Hyperlinks.Add($ws3.Range("A2"), ($srcBase + $file2), "", "", $file2)
$ws3.Hyperlinks.Add($ws3.Range("B2"), ($srcBase + $file2), "", "", ".png")
$ws3.Hyperlinks.Add($ws3.Range("D2"), ($deHtBase + $deTarget2), "", "", $deTarget2)

$ws3.Hyperlinks.Add($ws3.Range("A3"), ($srcBase + $file3), "", "", $file3)
$ws3.Hyperlinks.Add($ws3.Range("B3"), ($srcBase + $file3), "", "", ".png")
$ws3.Hyperlinks.Add($ws3.Range("D3"), ($deHtBase + $deTarget3), "", "", $deTarget3)

$ws3.Hyperlinks.Add($ws3.Range("A4"), ($srcBase + $file4), "", "", $file4)
$ws3.Hyperlinks.Add($ws3.Range("B4"), ($srcBase + $file4), "", "", ".md")
$ws3.Hyperlinks.Add($ws3.Range("D4"), ($deHtBase + $deTarget4), "", "", $deTarget4)
